$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("disk_list")

$ws.Range("A3").Value = "bez jmena a hesla"
$ws.Range("B3").Value = "P"
$ws.Range("C3").Value = "\\192.168.000.000\"
